$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 "invalidCredentialTest": the test user in row 2 changes from
# "john" to "bala".
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("invalidCredentialTest")
$ws1.Range("A2").Value = "bala"

# ---------------------------------------------------------------------
# Sheet2: rename "patientErrorMessageTest" -> "invalidCredentialTestCount"
# and populate it with a summary row + the same credentials table plus
# two extra rows (paul / mark).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("patientErrorMessageTest")
$ws2.Name = "invalidCredentialTestCount"

$ws2.Range("A1").Value = "Rows to be tested"
$ws2.Range("B1").Value = 2

$ws2.Range("A2").Value = "Username"
$ws2.Range("B2").Value = "Password"
$ws2.Range("C2").Value = "Language"
$ws2.Range("D2").Value = "ExpectedValue"

$ws2.Range("A3").Value = "bala"
$ws2.Range("B3").Value = "john123"
$ws2.Range("C3").Value = "Dutch"
$ws2.Range("D3").Value = "Invalid username or password"

$ws2.Range("A4").Value = "peter"
$ws2.Range("B4").Value = "peter123"
$ws2.Range("C4").Value = "Danish"
$ws2.Range("D4").Value = "Invalid username or password"

$ws2.Range("A5").Value = "paul"
$ws2.Range("B5").Value = "peter123"
$ws2.Range("C5").Value = "Danish"
$ws2.Range("D5").Value = "Invalid username or password"

$ws2.Range("A6").Value = "mark"
$ws2.Range("B6").Value = "peter123"
$ws2.Range("C6").Value = "Danish"
$ws2.Range("D6").Value = "Invalid username or password"

# Column widths mirroring the best-fit widths Excel computed for this
# content (A/B/D were auto-sized, C was set to a fixed custom width).
$ws2.Columns.Item(1).ColumnWidth = 16.1666666666667
$ws2.Columns.Item(2).ColumnWidth = 8.66666666666667
$ws2.Columns.Item(3).ColumnWidth = 11.5
$ws2.Columns.Item(4).ColumnWidth = 27.1666666666667

# ---------------------------------------------------------------------
# Selections / active sheet: sheet2 becomes the active tab with A2
# selected; sheet1 loses its tab selection and now has the whole table
# A1:D3 selected instead of the old B2 single-cell selection.
# ---------------------------------------------------------------------
$ws1.Range("A1:D3").Select()
$ws2.Activate()
$ws2.Range("A2").Select()
